$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the worksheet's used range
$lastRow = $ws.UsedRange.Rows.Count

# Column C holds the "Förändrad" (changed) date, stored as serial 46061 (2026-02-08).
# Increment every populated cell in column C (rows 2..lastRow) by one day to 46062 (2026-02-09).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
